$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "SONDAGEM MISTA"
$ws.Cells.Item(1, 2).Value = "MEDIÇÃO"
$ws.Cells.Item(1, 3).Value = "NORTE"
$ws.Cells.Item(1, 4).Value = "LESTE"

# Copy the header style from A1 (bold/border/center) across the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)

# --- Clear out now-unused columns E:M for rows 1-5 and rows 4-5 entirely ---
$ws.Range("E1:M5").Clear()
$ws.Range("A4:D5").Clear()

# --- Row 2 data ---
$ws.Cells.Item(2, 1).Value = "SM-km050-001"
$ws.Cells.Item(2, 2).Value = 45658
$ws.Cells.Item(2, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 3).Value = "N=250362.5494"
$ws.Cells.Item(2, 4).Value = "E=151188.7322"

# --- Row 3 data ---
$ws.Cells.Item(3, 1).Value = "SM-km050-002"
$ws.Cells.Item(3, 2).Value = 45658
$ws.Cells.Item(3, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 3).Value = "N=250248.7755"
$ws.Cells.Item(3, 4).Value = "E=151378.0133"

$wb.Save()
